$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Row 16: merge the Additional Effort [h] value into Effort [h] and clear the additional effort cell
$ws.Range("B16").Value = 5.5
$ws.Range("C16").ClearContents()

# New row 32: date, effort hours, and description referencing a new shared string
$ws.Range("A31").Copy($ws.Range("A32"))
$ws.Range("A32").Value = [DateTime]::FromOADate(41208)
$ws.Range("B32").Value = 4
$ws.Range("D32").Value = "Manual continued, widely completed as draft but without chapter Use Cases"
